$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" everywhere it appears ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"

# --- Narrow the (now shorter) status columns ---
# Overview sheet: columns E and F
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
